$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.322.66"
$ws.Range("D3").Value = "1.665.17"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.96%  "
$ws.Range("D5").Value = "'219.17"
$ws.Range("D6").Value = "'0.5348"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").Value = "'0.2659"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("D9").Value = "'0.06414"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'20.73"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").Value = "'0.07850"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "'4.565"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "1.665.11"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "1.892.67"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'0.5534"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "0.0₅8192"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'65.82"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "26.343.82"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "'4.686"
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "'193.33"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").Value = "'6.039"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'7.214"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "'1.498"
$ws.Range("E29").Value = "  +4.94%  "
$ws.Range("D30").Value = "'0.05872"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "'1.283"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").Value = "'3.634"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D34").Value = "'1.608"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").Value = "'0.9693"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").Value = "'2.828"
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("D37").Value = "'2.421"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "'0.5834"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'0.8710"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").Value = "'5.849"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.053.80"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'105.21"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").Value = "1.804.03"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").Value = "'57.79"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  -4.81%  "
$ws.Range("D48").Value = "'1.015"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "'0.4387"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "'8.005"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  +0.46%  "
